$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Create angular service") is now marked as Done:
# copy the formatting used by the other "Done" rows (e.g. row 9) onto row 6,
# and set the Status cell to "Done".
$ws.Range("A9:C9").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C6").Value = "Done"

# Row 8 ("Call server service from angular") is no longer Done:
# copy the formatting used by the "in progress" rows (e.g. row 7) onto row 8,
# and clear the Status cell.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C8").ClearContents()

$excel.CutCopyMode = $false

# Update the active selection to G8
$ws.Range("G8").Select()
